$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row=2; Col="D"; Value="29.189.67"; ForceText=$false }
    @{ Row=2; Col="E"; Value="  -0.79%  "; ForceText=$false }
    @{ Row=3; Col="D"; Value="1.860.50"; ForceText=$false }
    @{ Row=3; Col="E"; Value="  -0.82%  "; ForceText=$false }
    @{ Row=4; Col="D"; Value="0.9995"; ForceText=$true }
    @{ Row=4; Col="E"; Value="  -0.10%  "; ForceText=$false }
    @{ Row=5; Col="D"; Value="0.7082"; ForceText=$true }
    @{ Row=5; Col="E"; Value="  -0.57%  "; ForceText=$false }
    @{ Row=6; Col="D"; Value="240.88"; ForceText=$true }
    @{ Row=6; Col="E"; Value="  -0.38%  "; ForceText=$false }
    @{ Row=8; Col="D"; Value="0.3080"; ForceText=$true }
    @{ Row=8; Col="E"; Value="  -0.90%  "; ForceText=$false }
    @{ Row=9; Col="D"; Value="0.07650"; ForceText=$true }
    @{ Row=9; Col="E"; Value="  -2.51%  "; ForceText=$false }
    @{ Row=10; Col="D"; Value="24.75"; ForceText=$true }
    @{ Row=10; Col="E"; Value="  -1.65%  "; ForceText=$false }
    @{ Row=11; Col="D"; Value="0.08417"; ForceText=$true }
    @{ Row=11; Col="E"; Value="  +2.13%  "; ForceText=$false }
    @{ Row=12; Col="D"; Value="1.882.17"; ForceText=$false }
    @{ Row=12; Col="E"; Value="  +1.42%  "; ForceText=$false }
    @{ Row=13; Col="D"; Value="5.187"; ForceText=$true }
    @{ Row=13; Col="E"; Value="  -1.36%  "; ForceText=$false }
    @{ Row=14; Col="D"; Value="0.7095"; ForceText=$true }
    @{ Row=14; Col="E"; Value="  -2.45%  "; ForceText=$false }
    @{ Row=15; Col="D"; Value="91.18"; ForceText=$true }
    @{ Row=15; Col="E"; Value="  +0.45%  "; ForceText=$false }
    @{ Row=16; Col="D"; Value="29.199.28"; ForceText=$false }
    @{ Row=16; Col="E"; Value="  -0.78%  "; ForceText=$false }
    @{ Row=17; Col="D"; Value="5.924"; ForceText=$true }
    @{ Row=17; Col="E"; Value="  +0.36%  "; ForceText=$false }
    @{ Row=18; Col="D"; Value="242.87"; ForceText=$true }
    @{ Row=18; Col="E"; Value="  -1.77%  "; ForceText=$false }
    @{ Row=19; Col="D"; Value="0.000007812"; ForceText=$true }
    @{ Row=19; Col="E"; Value="  -0.60%  "; ForceText=$false }
    @{ Row=20; Col="D"; Value="2.114.64"; ForceText=$false }
    @{ Row=20; Col="E"; Value="  +0.25%  "; ForceText=$false }
    @{ Row=21; Col="E"; Value="  -1.31%  "; ForceText=$false }
    @{ Row=22; Col="E"; Value="  -0.03%  "; ForceText=$false }
    @{ Row=23; Col="D"; Value="7.870"; ForceText=$true }
    @{ Row=23; Col="E"; Value="  -1.15%  "; ForceText=$false }
    @{ Row=24; Col="D"; Value="0.9997"; ForceText=$true }
    @{ Row=24; Col="E"; Value="  -0.06%  "; ForceText=$false }
    @{ Row=25; Col="D"; Value="0.1591"; ForceText=$true }
    @{ Row=25; Col="E"; Value="  +1.07%  "; ForceText=$false }
    @{ Row=26; Col="D"; Value="163.11"; ForceText=$true }
    @{ Row=26; Col="E"; Value="  -0.26%  "; ForceText=$false }
    @{ Row=27; Col="D"; Value="8.919"; ForceText=$true }
    @{ Row=27; Col="E"; Value="  -0.72%  "; ForceText=$false }
    @{ Row=28; Col="D"; Value="18.43"; ForceText=$true }
    @{ Row=28; Col="E"; Value="  +0.86%  "; ForceText=$false }
    @{ Row=29; Col="E"; Value="  +0.47%  "; ForceText=$false }
    @{ Row=30; Col="D"; Value="1.316"; ForceText=$true }
    @{ Row=30; Col="E"; Value="  -3.54%  "; ForceText=$false }
    @{ Row=31; Col="D"; Value="4.398"; ForceText=$true }
    @{ Row=31; Col="E"; Value="  +0.91%  "; ForceText=$false }
    @{ Row=32; Col="D"; Value="4.219"; ForceText=$true }
    @{ Row=32; Col="E"; Value="  +2.34%  "; ForceText=$false }
    @{ Row=33; Col="D"; Value="0.05124"; ForceText=$true }
    @{ Row=33; Col="E"; Value="  -3.32%  "; ForceText=$false }
    @{ Row=34; Col="D"; Value="0.8060"; ForceText=$true }
    @{ Row=34; Col="E"; Value="  +11.77%  "; ForceText=$false }
    @{ Row=35; Col="D"; Value="1.905"; ForceText=$true }
    @{ Row=35; Col="E"; Value="  -1.13%  "; ForceText=$false }
    @{ Row=36; Col="E"; Value="  -2.94%  "; ForceText=$false }
    @{ Row=37; Col="D"; Value="2.680"; ForceText=$true }
    @{ Row=37; Col="E"; Value="  +0.06%  "; ForceText=$false }
    @{ Row=38; Col="E"; Value="  -1.01%  "; ForceText=$false }
    @{ Row=39; Col="E"; Value="  -0.88%  "; ForceText=$false }
    @{ Row=40; Col="D"; Value="1.167.06"; ForceText=$false }
    @{ Row=40; Col="E"; Value="  -6.58%  "; ForceText=$false }
    @{ Row=41; Col="D"; Value="6.182"; ForceText=$true }
    @{ Row=41; Col="E"; Value="  +0.65%  "; ForceText=$false }
    @{ Row=42; Col="D"; Value="0.8919"; ForceText=$true }
    @{ Row=42; Col="E"; Value="  -1.78%  "; ForceText=$false }
    @{ Row=43; Col="D"; Value="72.76"; ForceText=$true }
    @{ Row=43; Col="E"; Value="  -1.27%  "; ForceText=$false }
    @{ Row=44; Col="D"; Value="0.9997"; ForceText=$true }
    @{ Row=44; Col="E"; Value="  -0.11%  "; ForceText=$false }
    @{ Row=45; Col="D"; Value="101.93"; ForceText=$true }
    @{ Row=45; Col="E"; Value="  -1.23%  "; ForceText=$false }
    @{ Row=46; Col="D"; Value="2.013.13"; ForceText=$false }
    @{ Row=47; Col="D"; Value="0.5174"; ForceText=$true }
    @{ Row=47; Col="E"; Value="  -2.93%  "; ForceText=$false }
    @{ Row=48; Col="D"; Value="1.771"; ForceText=$true }
    @{ Row=48; Col="E"; Value="  +0.27%  "; ForceText=$false }
    @{ Row=49; Col="E"; Value="  +0.16%  "; ForceText=$false }
    @{ Row=50; Col="D"; Value="9.255"; ForceText=$true }
    @{ Row=50; Col="E"; Value="  +0.27%  "; ForceText=$false }
    @{ Row=51; Col="E"; Value="  +0.36%  "; ForceText=$false }
)

foreach ($change in $changes) {
    $cellRef = "$($change.Col)$($change.Row)"
    $cell = $ws.Range($cellRef)
    if ($change.ForceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $change.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $change.Value
    }
}
